$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grab the raw (unformatted) values and formats from the current last data row (89)
$dateVal = $ws.Range("A89").Value2
$bVal = $ws.Range("B89").Value2
$cVal = $ws.Range("C89").Value2
$dVal = $ws.Range("D89").Value2

$lastRowDateFormat = $ws.Range("A89").NumberFormat   # date-only format, used for the final row
$normalDateFormat = $ws.Range("A88").NumberFormat    # date+time format, used for all other rows

# Row 89 is no longer the last row: switch its date cell to the regular
# (non-final) date/time number format used by the rest of the data rows.
$ws.Range("A89").NumberFormat = $normalDateFormat

# Append new row 90 with the same data, which becomes the new last row.
$ws.Range("A90").Value2 = $dateVal
$ws.Range("B90").Value2 = $bVal
$ws.Range("C90").Value2 = $cVal
$ws.Range("D90").Value2 = $dVal

# New last row's date cell gets the "final row" date-only number format.
$ws.Range("A90").NumberFormat = $lastRowDateFormat
